$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.30321741104126
$ws.Range("B1").Value = 2.438845872879028
$ws.Range("C1").Value = 4.559343814849854
$ws.Range("D1").Value = 1.917994976043701
$ws.Range("E1").Value = 1.114999532699585
